# إضافة حدث جديد في Card23 by admin at 2025-12-08 07:30:59
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card23")

# Row 14 (existing row) currently has blank placeholder cells B14:K14.
# Fill them with the literal text "nan" to match the rest of the sheet's
# "missing value" convention used throughout the table.
$cols = @("B", "C", "D", "E", "F", "G", "H", "I", "J", "K")
foreach ($col in $cols) {
    $ws.Range($col + "14").Value = "nan"
}

# Row 15: a brand-new service event appended below row 14.
# Column A holds the card number as text ("23"); prefix with an
# apostrophe so Excel stores it as text instead of coercing to a number
# (matches how the other "card" column cells are stored).
$ws.Range("A15").Value = "'23"
$ws.Range("A15").ClearFormats()

# B15:K15 stay blank (no measurement data recorded for this event) -
# touch each so the cell is materialised in the sheet (matches the
# template's convention of emitting an empty cell per column).
foreach ($col in $cols) {
    $ws.Range($col + "15").Font.Bold = $false
}
$ws.Range("L15").Value = "6\8\2024"
# M15 (Event/hours) stays blank for this entry.
$ws.Range("M15").Font.Bold = $false
$ws.Range("N15").Value = "تم عمل setting كامل للمكنه وتضيق المسافات"
$ws.Range("O15").Value = "الخبير"
